$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(184).Insert()

$ws.Cells.Item(184, 1).Value = 5
$ws.Cells.Item(184, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(184, 3).Value = "Maule"
$ws.Cells.Item(184, 4).Value = 44726
$ws.Cells.Item(184, 5).Value = 7
$ws.Cells.Item(184, 6).Value = 100112009
$ws.Cells.Item(184, 7).Value = "Acelga"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 500
$ws.Cells.Item(184, 11).Value = 2500
$ws.Cells.Item(184, 12).Value = 2500
$ws.Cells.Item(184, 13).Value = 2500
$ws.Cells.Item(184, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(184, 15).Value = "Región del Maule"
$ws.Cells.Item(184, 16).Value = 625
$ws.Cells.Item(184, 17).Value = 4
$ws.Cells.Item(184, 18).Value = "Hortaliza"
